$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("G2").Value = -1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = -1
$ws.Range("J2").Value = 8
